$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 17 ("Assurance MM1 DCA" milestone field),
# pushing the existing "Project MM18 ..." block (previously rows 17-25)
# down to rows 18-26.
$ws.Rows("17:17").Insert()

# Label the new row in column A. B17/C17/D17 keep the formatting that was
# cloned from the row above on insert (styles 1/2/4), and are left blank.
$ws.Range("A17").Value = "Assurance MM1 DCA"

# The new row only has data in A:D - make sure E17/F17 carry no leftover
# formatting/content so they don't exist as cells at all.
$ws.Range("E17:F17").Clear()

# Update the remembered selection on the bottom-right frozen pane.
$null = $ws.Range("K24").Select()
